# Fix calculation errors in the "Final Costing" section of the Costing sheet.
# The old "Monthly Cost" formula accidentally multiplied the output-token term by
# L25 (Input Length) a second time and folded "articles/day * 31" directly into a
# single cell, which was off by a factor of ~30. This rewrites the section as a
# clear chain: Cost / Article -> Cost / Day -> Cost / Month, plus a GBP
# conversion column using a "Dollars to pounds" rate.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the TL;DR note cells (row 24/25/27 unchanged text, row 28 note removed) ---
$ws.Range("D28").ClearContents()

# --- Update the input parameters ---
$ws.Range("L24").Value = 15      # Articles / Day: 5 -> 15
$ws.Range("L25").Value = 3500    # Input Length: 1500 -> 3500
$ws.Range("L27").Value = 100     # Output Length: 50 -> 100
# L28 (Safety Factor = 1.2) is unchanged

# --- Clear the old "Monthly Cost" row (row 30) and its special thick-bottom formatting ---
$ws.Range("K29:L30").Clear()

# --- New row 29: Days / Month Active ---
$ws.Range("K29").Value = "Days / Month Active"
$ws.Range("K29").HorizontalAlignment = -4108
$ws.Range("L29").Value = 31
$ws.Range("K29:L29").Borders.LineStyle = 1
$ws.Range("K29:L29").Borders.Weight = 2
$ws.Range("K29:L29").VerticalAlignment = -4108

# --- Row 31: Cost / Article ---
$ws.Range("K31").Value = "Cost / Article"
$ws.Range("L31").Formula = "=L28*(L25*VLOOKUP(L26,C4:E6,2,FALSE)+VLOOKUP(L26,C4:E6,3,FALSE)*L27)*D8/1000000"
$ws.Range("L31").NumberFormat = '_-[$$-409]* #,##0.0000_ ;_-[$$-409]* \-#,##0.0000\ ;_-[$$-409]* "-"??_ ;_-@_ '
$ws.Range("M31").Formula = "=L31*`$L`$36"

# --- Row 32: Cost / Day ---
$ws.Range("K32").Value = "Cost / Day"
$ws.Range("L32").Formula = "=L31*L24"
$ws.Range("L32").NumberFormat = '_-[$$-409]* #,##0.0000_ ;_-[$$-409]* \-#,##0.0000\ ;_-[$$-409]* "-"????_ ;_-@_ '
$ws.Range("M32").Formula = "=L32*`$L`$36"

# --- Row 33: Cost / Month ---
$ws.Range("K33").Value = "Cost / Month"
$ws.Range("L33").Formula = "=L32*L29"
$ws.Range("L33").NumberFormat = '_-[$$-409]* #,##0.0000_ ;_-[$$-409]* \-#,##0.0000\ ;_-[$$-409]* "-"????_ ;_-@_ '
$ws.Range("M33").Formula = "=L33*`$L`$36"

# --- Style K31:L33 like the rest of the input table (thin box border, normal weight) ---
$ws.Range("K31:L33").Borders.LineStyle = 1
$ws.Range("K31:L33").Borders.Weight = 2
$ws.Range("K31:L33").VerticalAlignment = -4108

# --- GBP column (M31:M33): bold, euro-bracket style currency format used for the converted figures ---
$ws.Range("M31:M33").NumberFormat = '_-* #,##0.0000\ [$€-1]_-;\-* #,##0.0000\ [$€-1]_-;_-* "-"??\ [$€-1]_-;_-@_-'
$ws.Range("M31:M33").Font.Bold = $true
$ws.Range("M31:M33").Borders.LineStyle = 1
$ws.Range("M31:M33").Borders.Weight = 2
$ws.Range("M31:M33").VerticalAlignment = -4108

# --- Row 36: Dollars to pounds conversion rate ---
$ws.Range("K36").Value = "Dollars to pounds"
$ws.Range("L36").Value = 0.79

$wb.Save()
